$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Turn the " m:1/0 " field ( { m:1/0 } ) into plain literal text
#    "{" "m" ":1/0" "}" split across four runs (no field codes left).
# ---------------------------------------------------------------------------
$f = $d.Fields.Item(1)
$code = $f.Code
$fieldPos = $code.Start

# Inserting text at any position inside the field's code range places the
# new text immediately in front of the field (outside of it), as its own
# unformatted run(s).
$ins = $d.Range($fieldPos, $fieldPos)
$ins.InsertBefore("{m:1/0}")

# The text above landed right before the field, i.e. at the position where
# the field used to begin (use the ORIGINAL position captured before the
# insertion shifted everything after it).
$base = $fieldPos - 1

$r1 = $d.Range($base, $base + 1)       # "{"
$r2 = $d.Range($base + 1, $base + 2)   # "m"
$r3 = $d.Range($base + 2, $base + 6)   # ":1/0"
$r4 = $d.Range($base + 6, $base + 7)   # "}"

# Force run boundaries between the four pieces (otherwise Word silently
# merges adjacent same-format runs into a single run), then put the
# formatting back the way it started so no visible formatting changes.
$r2.Font.Bold = 1
$r4.Font.Bold = 1
$r2.Font.Bold = 0
$r4.Font.Bold = 0

# Remove the now orphaned field (begin/instrText/end runs).
$f.Delete()

# ---------------------------------------------------------------------------
# 2. Prepend "    <---" to the error message text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "    <---divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Refresh the stale stack-trace line numbers / accessor names.
# ---------------------------------------------------------------------------
function Replace-Literal($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Literal "java.lang.reflect.Method.invoke(Method.java:568)" "java.lang.reflect.Method.invoke(Method.java:569)"
Replace-Literal "M2DocEvaluator.caseQuery(M2DocEvaluator.java:604)" "M2DocEvaluator.caseQuery(M2DocEvaluator.java:659)"
Replace-Literal "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1459)" "M2DocEvaluator.doSwitch(M2DocEvaluator.java:2022)"
Replace-Literal "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1684)" "M2DocEvaluator.caseBlock(M2DocEvaluator.java:2247)"
Replace-Literal "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:314)" "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:340)"
Replace-Literal "M2DocEvaluator.generate(M2DocEvaluator.java:299)" "M2DocEvaluator.generate(M2DocEvaluator.java:324)"
Replace-Literal "M2DocUtils.generate(M2DocUtils.java:853)" "M2DocUtils.generate(M2DocUtils.java:912)"
Replace-Literal "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:508)" "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:497)"
Replace-Literal "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:400)" "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:391)"
Replace-Literal "jdk.internal.reflect.GeneratedMethodAccessor6.invoke" "jdk.internal.reflect.GeneratedMethodAccessor7.invoke"
Replace-Literal "RemoteTestRunner.runTests(RemoteTestRunner.java:756)" "RemoteTestRunner.runTests(RemoteTestRunner.java:757)"

Write-Output "done"
